$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (closest achievable via ColumnWidth quantization)
$ws.Columns.Item(5).ColumnWidth = 21.833333333333332
$ws.Columns.Item(7).ColumnWidth = 20.833333333333332
$ws.Columns.Item(9).ColumnWidth = 19.833333333333332
$ws.Columns.Item(15).ColumnWidth = 20.833333333333332

# Data value updates rows 2-12, columns A-T (recalculated ratios after adding laboratory numbers)
$ws.Cells.Item(2, 1).Value2 = -105.6004059644875
$ws.Cells.Item(2, 2).Value2 = 0.001522748921522226
$ws.Cells.Item(2, 3).Value2 = 0.01023542238374481
$ws.Cells.Item(2, 4).Value2 = 2.525366425064601
$ws.Cells.Item(2, 5).Value2 = 0.007165735470860413
$ws.Cells.Item(2, 6).Value2 = 0.01478863668832112
$ws.Cells.Item(2, 7).Value2 = 2.732165266071587
$ws.Cells.Item(2, 8).Value2 = 0.02316925422885041
$ws.Cells.Item(2, 9).Value2 = 1.808951687183812
$ws.Cells.Item(2, 10).Value2 = 2.576467054448028
$ws.Cells.Item(2, 11).Value2 = 0.006778780467003105
$ws.Cells.Item(2, 12).Value2 = 0.1702537581274621
$ws.Cells.Item(2, 13).Value2 = 0.00004916399262409689
$ws.Cells.Item(2, 14).Value2 = 0.1702537581274623
$ws.Cells.Item(2, 15).Value2 = 1.318240643972384
$ws.Cells.Item(2, 16).Value2 = 0.5337376544124303
$ws.Cells.Item(2, 17).Value2 = 0.1702167991392523
$ws.Cells.Item(2, 18).Value2 = 1.613088256207097
$ws.Cells.Item(2, 19).Value2 = 0.2230016042121955
$ws.Cells.Item(2, 20).Value2 = 1.55796235410836
$ws.Cells.Item(3, 1).Value2 = 0.1774019644116276
$ws.Cells.Item(3, 2).Value2 = 0.002107527865453386
$ws.Cells.Item(3, 3).Value2 = 0.01023651243346234
$ws.Cells.Item(3, 4).Value2 = 3.9025283601347
$ws.Cells.Item(3, 5).Value2 = 0.007151741669928597
$ws.Cells.Item(3, 6).Value2 = 0.243603262622918
$ws.Cells.Item(3, 7).Value2 = 0.2056506777876824
$ws.Cells.Item(3, 8).Value2 = 2.826771272669783
$ws.Cells.Item(3, 9).Value2 = 0.1561259601382454
$ws.Cells.Item(3, 10).Value2 = 2.240823406699239
$ws.Cells.Item(3, 11).Value2 = 0.007580485368271605
$ws.Cells.Item(3, 12).Value2 = 0.2107154052185198
$ws.Cells.Item(3, 13).Value2 = 0.0000549784623571892
$ws.Cells.Item(3, 14).Value2 = 0.2107154052185211
$ws.Cells.Item(3, 15).Value2 = 0.03604320149520462
$ws.Cells.Item(3, 16).Value2 = 5.579804477398916
$ws.Cells.Item(3, 17).Value2 = 0.327277029526253
$ws.Cells.Item(3, 18).Value2 = 1.221362702752062
$ws.Cells.Item(3, 19).Value2 = 0.01154992788514052
$ws.Cells.Item(3, 20).Value2 = 0.779407971311902
$ws.Cells.Item(4, 1).Value2 = -104.3468227271361
$ws.Cells.Item(4, 2).Value2 = 0.001189460449408623
$ws.Cells.Item(4, 3).Value2 = 0.009931202987964183
$ws.Cells.Item(4, 4).Value2 = 2.847590900310127
$ws.Cells.Item(4, 5).Value2 = 0.007161373392654631
$ws.Cells.Item(4, 6).Value2 = 0.01241336652264171
$ws.Cells.Item(4, 7).Value2 = 2.73261122757083
$ws.Cells.Item(4, 8).Value2 = 0.02835810671222438
$ws.Cells.Item(4, 9).Value2 = 1.868456204221887
$ws.Cells.Item(4, 10).Value2 = 3.005917901200326
$ws.Cells.Item(4, 11).Value2 = 0.006788281550880812
$ws.Cells.Item(4, 12).Value2 = 0.1328036878103157
$ws.Cells.Item(4, 13).Value2 = 0.00004923290047853448
$ws.Cells.Item(4, 14).Value2 = 0.1328036878103152
$ws.Cells.Item(4, 15).Value2 = 1.31762407900823
$ws.Cells.Item(4, 16).Value2 = 0.7002340683123393
$ws.Cells.Item(4, 17).Value2 = 0.1724751551947613
$ws.Cells.Item(4, 18).Value2 = 1.285509659675717
$ws.Cells.Item(4, 19).Value2 = 0.227718299854116
$ws.Cells.Item(4, 20).Value2 = 1.001689343495457
$ws.Cells.Item(5, 1).Value2 = -7.452253930086417
$ws.Cells.Item(5, 2).Value2 = 0.004968873446402617
$ws.Cells.Item(5, 3).Value2 = 0.01003286248391662
$ws.Cells.Item(5, 4).Value2 = 0.4349841647427725
$ws.Cells.Item(5, 5).Value2 = 0.00713951470888153
$ws.Cells.Item(5, 6).Value2 = 0.124516000850981
$ws.Cells.Item(5, 7).Value2 = 0.03411091781235923
$ws.Cells.Item(5, 8).Value2 = 0.111622494396719
$ws.Cells.Item(5, 9).Value2 = 0.02564128291718566
$ws.Cells.Item(5, 10).Value2 = 0.659620017780448
$ws.Cells.Item(5, 11).Value2 = 0.007522659131886342
$ws.Cells.Item(5, 12).Value2 = 0.5006180776771035
$ws.Cells.Item(5, 13).Value2 = 0.00005455907000882169
$ws.Cells.Item(5, 14).Value2 = 0.5006180776771032
$ws.Cells.Item(5, 15).Value2 = 0.00581951220202368
$ws.Cells.Item(5, 16).Value2 = 2.780815804282748
$ws.Cells.Item(5, 17).Value2 = 0.224860292588127
$ws.Cells.Item(5, 18).Value2 = 1.324936855350163
$ws.Cells.Item(5, 19).Value2 = 0.001290716398291085
$ws.Cells.Item(5, 20).Value2 = 2.940693528015276
$ws.Cells.Item(6, 1).Value2 = -62.455539026395
$ws.Cells.Item(6, 2).Value2 = 0.001339707136559375
$ws.Cells.Item(6, 3).Value2 = 0.009938835616045487
$ws.Cells.Item(6, 4).Value2 = 2.816784145733791
$ws.Cells.Item(6, 5).Value2 = 0.007138558533294139
$ws.Cells.Item(6, 6).Value2 = 0.01676722484027121
$ws.Cells.Item(6, 7).Value2 = 2.732135489862179
$ws.Cells.Item(6, 8).Value2 = 0.03758562717671422
$ws.Cells.Item(6, 9).Value2 = 1.951132635891098
$ws.Cells.Item(6, 10).Value2 = 2.610342952444355
$ws.Cells.Item(6, 11).Value2 = 0.007105781488919688
$ws.Cells.Item(6, 12).Value2 = 0.1428953177503862
$ws.Cells.Item(6, 13).Value2 = 0.00005153561033731759
$ws.Cells.Item(6, 14).Value2 = 0.1428953177503859
$ws.Cells.Item(6, 15).Value2 = 1.315685602038329
$ws.Cells.Item(6, 16).Value2 = 0.5409586722616812
$ws.Cells.Item(6, 17).Value2 = 0.1815974558276844
$ws.Cells.Item(6, 18).Value2 = 1.27636062647053
$ws.Cells.Item(6, 19).Value2 = 0.2399428959188615
$ws.Cells.Item(6, 20).Value2 = 1.279373571369847
$ws.Cells.Item(7, 1).Value2 = 8.558799757307334
$ws.Cells.Item(7, 2).Value2 = 0.001978891721470877
$ws.Cells.Item(7, 3).Value2 = 0.01000159070118261
$ws.Cells.Item(7, 4).Value2 = 0.2995195323465933
$ws.Cells.Item(7, 5).Value2 = 0.00714736006040352
$ws.Cells.Item(7, 6).Value2 = 0.03217463243436037
$ws.Cells.Item(7, 7).Value2 = 0.1302778757219377
$ws.Cells.Item(7, 8).Value2 = 0.03044843328979099
$ws.Cells.Item(7, 9).Value2 = 0.09969720717163566
$ws.Cells.Item(7, 10).Value2 = 0.4051568008543506
$ws.Cells.Item(7, 11).Value2 = 0.007644009162360457
$ws.Cells.Item(7, 12).Value2 = 0.1962098513192351
$ws.Cells.Item(7, 13).Value2 = 0.00005543917698856591
$ws.Cells.Item(7, 14).Value2 = 0.1962098513192353
$ws.Cells.Item(7, 15).Value2 = 0.01782589513152004
$ws.Cells.Item(7, 16).Value2 = 1.002151359149723
$ws.Cells.Item(7, 17).Value2 = 0.2501805126752944
$ws.Cells.Item(7, 18).Value2 = 0.5595788267333994
$ws.Cells.Item(7, 19).Value2 = 0.004467168957753317
$ws.Cells.Item(7, 20).Value2 = 0.6938224478811204
$ws.Cells.Item(8, 1).Value2 = -75.637059503742
$ws.Cells.Item(8, 2).Value2 = 0.001344655020757129
$ws.Cells.Item(8, 3).Value2 = 0.01003148007266415
$ws.Cells.Item(8, 4).Value2 = 1.995271778797226
$ws.Cells.Item(8, 5).Value2 = 0.007146590616651994
$ws.Cells.Item(8, 6).Value2 = 0.01366064184467083
$ws.Cells.Item(8, 7).Value2 = 2.731930548347546
$ws.Cells.Item(8, 8).Value2 = 0.03047594804747604
$ws.Cells.Item(8, 9).Value2 = 1.91296454162027
$ws.Cells.Item(8, 10).Value2 = 2.054691392795934
$ws.Cells.Item(8, 11).Value2 = 0.007005876889081848
$ws.Cells.Item(8, 12).Value2 = 0.145468296255499
$ws.Cells.Item(8, 13).Value2 = 0.00005081103915029519
$ws.Cells.Item(8, 14).Value2 = 0.145468296255498
$ws.Cells.Item(8, 15).Value2 = 1.31376371035002
$ws.Cells.Item(8, 16).Value2 = 0.40047877722829
$ws.Cells.Item(8, 17).Value2 = 0.1796494562131981
$ws.Cells.Item(8, 18).Value2 = 1.36647001859519
$ws.Cells.Item(8, 19).Value2 = 0.2357733899885544
$ws.Cells.Item(8, 20).Value2 = 1.309927692333007
$ws.Cells.Item(9, 1).Value2 = 12.65834557701262
$ws.Cells.Item(9, 2).Value2 = 0.001617865900988797
$ws.Cells.Item(9, 3).Value2 = 0.01003530148248776
$ws.Cells.Item(9, 4).Value2 = 0.220048417453029
$ws.Cells.Item(9, 5).Value2 = 0.007143863131307234
$ws.Cells.Item(9, 6).Value2 = 0.01362275704189219
$ws.Cells.Item(9, 7).Value2 = 0.38626528275149
$ws.Cells.Item(9, 8).Value2 = 0.008381737472966153
$ws.Cells.Item(9, 9).Value2 = 0.2949953497862585
$ws.Cells.Item(9, 10).Value2 = 0.2453161363518621
$ws.Cells.Item(9, 11).Value2 = 0.007675080197400639
$ws.Cells.Item(9, 12).Value2 = 0.1597642391488841
$ws.Cells.Item(9, 13).Value2 = 0.00005566452373714028
$ws.Cells.Item(9, 14).Value2 = 0.1597642391488846
$ws.Cells.Item(9, 15).Value2 = 0.1441234321365172
$ws.Cells.Item(9, 16).Value2 = 0.3540798507298535
$ws.Cells.Item(9, 17).Value2 = 0.01852905203638091
$ws.Cells.Item(9, 18).Value2 = 0.2680031094358345
$ws.Cells.Item(9, 19).Value2 = 0.002669231605397931
$ws.Cells.Item(9, 20).Value2 = 0.2998133159033358
$ws.Cells.Item(10, 1).Value2 = -71.04536126314554
$ws.Cells.Item(10, 2).Value2 = 0.001324340229563539
$ws.Cells.Item(10, 3).Value2 = 0.009688541303555528
$ws.Cells.Item(10, 4).Value2 = 2.957211552912798
$ws.Cells.Item(10, 5).Value2 = 0.007142427503754504
$ws.Cells.Item(10, 6).Value2 = 0.01643797653779869
$ws.Cells.Item(10, 7).Value2 = 2.731666440643729
$ws.Cells.Item(10, 8).Value2 = 0.02987075204664381
$ws.Cells.Item(10, 9).Value2 = 1.998339578541245
$ws.Cells.Item(10, 10).Value2 = 2.966860945576796
$ws.Cells.Item(10, 11).Value2 = 0.007040678016621821
$ws.Cells.Item(10, 12).Value2 = 0.1425624217092339
$ws.Cells.Item(10, 13).Value2 = 0.000051063438882963
$ws.Cells.Item(10, 14).Value2 = 0.1425624217092355
$ws.Cells.Item(10, 15).Value2 = 1.315343398344421
$ws.Cells.Item(10, 16).Value2 = 0.5602296882894516
$ws.Cells.Item(10, 17).Value2 = 0.177974970287479
$ws.Cells.Item(10, 18).Value2 = 1.348432371958944
$ws.Cells.Item(10, 19).Value2 = 0.2340376254952511
$ws.Cells.Item(10, 20).Value2 = 1.249778956658016
$ws.Cells.Item(11, 1).Value2 = 83.89337830856735
$ws.Cells.Item(11, 2).Value2 = 0.001280124441955429
$ws.Cells.Item(11, 3).Value2 = 0.009530949851512161
$ws.Cells.Item(11, 4).Value2 = 0.3029886373501832
$ws.Cells.Item(11, 5).Value2 = 0.00713321812664829
$ws.Cells.Item(11, 6).Value2 = 0.01409443544888798
$ws.Cells.Item(11, 7).Value2 = 0.449309322731185
$ws.Cells.Item(11, 8).Value2 = 0.009845090719645582
$ws.Cells.Item(11, 9).Value2 = 0.3876068192742979
$ws.Cells.Item(11, 10).Value2 = 0.3359403568749768
$ws.Cells.Item(11, 11).Value2 = 0.008214980541348931
$ws.Cells.Item(11, 12).Value2 = 0.1181042773739492
$ws.Cells.Item(11, 13).Value2 = 0.00005958022165018338
$ws.Cells.Item(11, 14).Value2 = 0.1181042773739502
$ws.Cells.Item(11, 15).Value2 = 0.1435802128450965
$ws.Cells.Item(11, 16).Value2 = 0.336225689862737
$ws.Cells.Item(11, 17).Value2 = 0.0009642275086030902
$ws.Cells.Item(11, 18).Value2 = 0.242283315585006
$ws.Cells.Item(11, 19).Value2 = 0.0001385718953360891
$ws.Cells.Item(11, 20).Value2 = 0.2439946920379288
$ws.Cells.Item(12, 1).Value2 = -55.71146207928979
$ws.Cells.Item(12, 2).Value2 = 0.001278282583173604
$ws.Cells.Item(12, 3).Value2 = 0.01012869368785669
$ws.Cells.Item(12, 4).Value2 = 2.810876262437888
$ws.Cells.Item(12, 5).Value2 = 0.007130956613605904
$ws.Cells.Item(12, 6).Value2 = 0.01752644314192083
$ws.Cells.Item(12, 7).Value2 = 2.733315956497029
$ws.Cells.Item(12, 8).Value2 = 0.03098775925022471
$ws.Cells.Item(12, 9).Value2 = 1.943086602579645
$ws.Cells.Item(12, 10).Value2 = 2.901203152336854
$ws.Cells.Item(12, 11).Value2 = 0.00715689579776092
$ws.Cells.Item(12, 12).Value2 = 0.1353699141565714
$ws.Cells.Item(12, 13).Value2 = 0.00005190632355263537
$ws.Cells.Item(12, 14).Value2 = 0.1353699141565709
$ws.Cells.Item(12, 15).Value2 = 1.31633189883816
$ws.Cells.Item(12, 16).Value2 = 0.4543621109275845
$ws.Cells.Item(12, 17).Value2 = 0.1786038163860654
$ws.Cells.Item(12, 18).Value2 = 1.294252663780376
$ws.Cells.Item(12, 19).Value2 = 0.2357903936137166
$ws.Cells.Item(12, 20).Value2 = 1.367741377955286
